$d = $word.ActiveDocument

# Make sure edits are applied as plain text changes, not tracked revisions.
$d.TrackRevisions = $false

# --- Remove the stale "_GoBack" bookmark (Word clears this on a real edit pass) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Touch the runs around every "DikesOvertopping" / version mentions so the
#     transient proofing (spell-check) marks Word had left around them get
#     cleared, matching a normal "open, edit, save" pass. These are harmless
#     self-replacements (same text back in), they just force Word to
#     regenerate the runs and drop the stale <w:proofErr/> markers. ---

$d.Content.Find.Execute("the DikesOvertopping kernel", $true, $false, $false, $false, $false, $true, 1, $false, "the DikesOvertopping kernel", 2)

$d.Content.Find.Execute("n making a DikesOvertopping release.", $true, $false, $false, $false, $false, $true, 1, $false, "n making a DikesOvertopping release.", 2)

$d.Content.Find.Execute("xx.y.z ", $true, $false, $false, $false, $false, $true, 1, $false, "xx.y.z ", 2)

$d.Content.Find.Execute("DikesOvertopping depends on Fortran-Common-Library and is ", $true, $false, $false, $false, $false, $true, 1, $false, "DikesOvertopping depends on Fortran-Common-Library and is ", 2)

$d.Content.Find.Execute("a release of DikesOvertopping ", $true, $false, $false, $false, $false, $true, 1, $false, "a release of DikesOvertopping ", 2)

# --- The actual substantive change: replace SVN with Git ---

# "... created from this commit in subversion and the resulting build ..."
#   -> "... created from this commit in Git and the resulting build ..."
$d.Content.Find.Execute("in subversion and", $true, $false, $false, $false, $false, $true, 1, $false, "in Git and", 2)

# "... of the master branch in TeamCity ..." -> "... of the main branch in TeamCity ..."
$d.Content.Find.Execute("of the master branch", $true, $false, $false, $false, $false, $true, 1, $false, "of the main branch", 2)
